# Update the workbook to the 2025-11-05 14:44 snapshot.
#
# 1) "Metadata" sheet: bump the "Last Updated" timestamp in A2.
# 2) "Stock List" sheet: a new instrument (CAPTRU-RE1) was inserted at the
#    top of the list (row 2). Every existing row shifts down by one, and
#    the table keeps its original size, so the previous last row
#    (TRAVELFOOD) drops off the bottom instead of growing the sheet.
#
# Note: this COM shim's `Range.Value` / `Cells.Item().Value` *getters* are
# not reliable (they surface a reflection placeholder instead of the real
# cell contents), so all reads below go through `.Value2`, which works
# correctly. Writes use the normal `.Value` setter.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata!A2 -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(2, 1).Value = "05 Nov 2025, 02:44 PM"

# --- 2) Stock List: insert CAPTRU-RE1 at row 2, shift the rest down -------
$ws = $wb.Worksheets.Item("Stock List")

$firstDataRow = 2
$lastDataRow = 76

# Snapshot the old values for rows 2..75 (these are the rows that will be
# shifted down into 3..76; the old row 76 is dropped).
$bVals = @()
$cVals = @()
$dVals = @()
$eVals = @()
$hVals = @()

for ($i = $firstDataRow; $i -le ($lastDataRow - 1); $i++) {
    $bVals += $ws.Cells.Item($i, 2).Value2
    $cVals += $ws.Cells.Item($i, 3).Value2
    $dVals += $ws.Cells.Item($i, 4).Value2
    $eVals += $ws.Cells.Item($i, 5).Value2
    $hVals += $ws.Cells.Item($i, 8).Value2
}

# Write back bottom-up so we never overwrite a row before it has been read.
for ($i = $lastDataRow; $i -ge ($firstDataRow + 1); $i--) {
    $srcIdx = $i - ($firstDataRow + 1)
    $ws.Cells.Item($i, 2).Value = $bVals[$srcIdx]
    $ws.Cells.Item($i, 3).Value = $cVals[$srcIdx]
    $ws.Cells.Item($i, 4).Value = $dVals[$srcIdx]
    $ws.Cells.Item($i, 5).Value = $eVals[$srcIdx]
    $ws.Cells.Item($i, 8).Value = $hVals[$srcIdx]
}

# New row at the top of the table.
$ws.Cells.Item($firstDataRow, 2).Value = "CAPTRU-RE1"
$ws.Cells.Item($firstDataRow, 3).Value = "CAPTRU-RE1"
$ws.Cells.Item($firstDataRow, 4).Value = 5.67
$ws.Cells.Item($firstDataRow, 5).Value = -11.9565
$ws.Cells.Item($firstDataRow, 8).Value = 0
